{"js": "// Apply the warrant-template wording fixes described by the diff:\n//  1. \"{{ YEARS_ON }}\" -> \"{{ YEARS }}\"\n//  2. Merge the \"Therefore, I request ... issued\" paragraph with the\n//     following \"{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTJUSTIFY }}\" paragraph\n//     into one paragraph, dropping the trailing colon and disambiguating\n//     the NIGHTTIME placeholder to NIGHTTIME1.\n//  3. \"Which offense occurred\" -> \"Which offense(s) occurred\"\n//  4. \"{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTKICKER }}\" ->\n//     \"{{ DAYTIME }}{{ NIGHTTIME2 }}\" (NIGHTKICKER placeholder dropped,\n//     NIGHTTIME disambiguated to NIGHTTIME2)\n\nconst body = context.document.body;\n\n// --- Change 1: {{ YEARS_ON }} -> {{ YEARS }} -------------------------------\nlet yearsResults = body.search(\"YEARS_ON\", { matchCase: true });\nyearsResults.load(\"items\");\nawait context.sync();\nif (yearsResults.items.length > 0) {\n  yearsResults.items[0].insertText(\"YEARS\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 2: merge the \"Therefore, I request...\" paragraph with the\n// following DAYTIME/NIGHTTIME/NIGHTJUSTIFY paragraph -----------------------\nlet issuedResults = body.search(\n  \"Therefore, I request that a search warrant be issued and a search be made for the item(s) described herein, and that the same be retained in my custody or the custody of the Oro Valley Police Department:\",\n  { matchCase: true }\n);\nissuedResults.load(\"items\");\nawait context.sync();\n\nif (issuedResults.items.length > 0) {\n  const firstPara = issuedResults.items[0].paragraphs.getFirst();\n  firstPara.load(\"text\");\n  await context.sync();\n\n  const nextPara = firstPara.getNext();\n  nextPara.load(\"text\");\n  await context.sync();\n\n  // nextPara.text is expected to be \"{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTJUSTIFY }}\"\n  const mergedDaytimeNighttime = nextPara.text.replace(\n    \"{{ NIGHTTIME }}\",\n    \"{{ NIGHTTIME1 }}\"\n  );\n\n  firstPara\n    .getRange()\n    .insertText(\n      \"Therefore, I request that a search warrant be issued and a search be made for the item(s) described herein, and that the same be retained in my custody or the custody of the Oro Valley Police Department \" +\n        mergedDaytimeNighttime,\n      Word.InsertLocation.replace\n    );\n  nextPara.getRange(\"Whole\").delete();\n  await context.sync();\n}\n\n// --- Change 3: \"Which offense occurred\" -> \"Which offense(s) occurred\" ----\nlet offenseResults = body.search(\"Which offense occurred \", { matchCase: true });\noffenseResults.load(\"items\");\nawait context.sync();\nif (offenseResults.items.length > 0) {\n  offenseResults.items[0].insertText(\n    \"Which offense(s) occurred \",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- Change 4: {{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTKICKER }} ->\n// {{ DAYTIME }}{{ NIGHTTIME2 }} ---------------------------------------------\nlet kickerResults = body.search(\"{{ NIGHTTIME }}{{ NIGHTKICKER }}\", {\n  matchCase: true,\n});\nkickerResults.load(\"items\");\nawait context.sync();\nif (kickerResults.items.length > 0) {\n  kickerResults.items[0].insertText(\n    \"{{ NIGHTTIME2 }}\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Apply the warrant-template wording fixes described by the diff:\n#  1. \"{{ YEARS_ON }}\" -> \"{{ YEARS }}\"\n#  2. Merge the \"Therefore, I request ... issued\" paragraph with the\n#     following \"{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTJUSTIFY }}\" paragraph\n#     into one paragraph, dropping the trailing colon and disambiguating\n#     the NIGHTTIME placeholder to NIGHTTIME1.\n#  3. \"Which offense occurred\" -> \"Which offense(s) occurred\"\n#  4. \"{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTKICKER }}\" ->\n#     \"{{ DAYTIME }}{{ NIGHTTIME2 }}\" (NIGHTKICKER placeholder dropped,\n#     NIGHTTIME disambiguated to NIGHTTIME2)\n\n$d = $word.ActiveDocument\n\n# --- Change 1: {{ YEARS_ON }} -> {{ YEARS }} --------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"YEARS_ON\", $false, $false, $false, $false, $false, $true, 1, $false, \"YEARS\", 2)\n\n# --- Change 2: merge the \"Therefore, I request...\" paragraph with the\n# following DAYTIME/NIGHTTIME/NIGHTJUSTIFY paragraph, dropping the colon\n# and the paragraph break, and renaming NIGHTTIME -> NIGHTTIME1 ------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"Department:^p{{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTJUSTIFY }}\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Department {{ DAYTIME }}{{ NIGHTTIME1 }}{{ NIGHTJUSTIFY }}\", 2)\n\n# --- Change 3: \"Which offense occurred\" -> \"Which offense(s) occurred\" ----\n$rng3 = $d.Content\n$rng3.Find.Execute(\"Which offense occurred \", $false, $false, $false, $false, $false, $true, 1, $false, \"Which offense(s) occurred \", 2)\n\n# --- Change 4: {{ DAYTIME }}{{ NIGHTTIME }}{{ NIGHTKICKER }} ->\n# {{ DAYTIME }}{{ NIGHTTIME2 }} (NIGHTKICKER dropped) -----------------------\n$rng4 = $d.Content\n$rng4.Find.Execute(\"{{ NIGHTTIME }}{{ NIGHTKICKER }}\", $false, $false, $false, $false, $false, $true, 1, $false, \"{{ NIGHTTIME2 }}\", 2)\n"}
